$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")  # row 28
$ws.Cells.Item(28, 8).Value = 591.7778
$ws.Cells.Item(28, 9).Value = 540
$ws.Cells.Item(28, 11).Value = 540
$ws.Cells.Item(28, 13).Value = -55

$ws = $wb.Worksheets.Item("ALC")  # row 86
$ws.Cells.Item(86, 8).Value = 0
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).Value = $null

$ws = $wb.Worksheets.Item("ALC")  # row 89
$ws.Cells.Item(89, 8).Value = 0
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).Value = $null

$ws = $wb.Worksheets.Item("ARM")  # row 125
$ws.Cells.Item(125, 8).Value = 94999.5
$ws.Cells.Item(125, 10).Value = 94999.5
$ws.Cells.Item(125, 12).Value = 94999.5
$ws.Cells.Item(125, 14).Value = -104839.5

$ws = $wb.Worksheets.Item("BSM")  # row 94
$ws.Cells.Item(94, 8).Value = 6801
$ws.Cells.Item(94, 9).Value = 6001.5
$ws.Cells.Item(94, 10).Value = 8400
$ws.Cells.Item(94, 11).Value = 6001.5
$ws.Cells.Item(94, 12).Value = 8400
$ws.Cells.Item(94, 13).Value = -5550.5
$ws.Cells.Item(94, 14).Value = -9302

$ws = $wb.Worksheets.Item("CRP")  # row 109
$ws.Cells.Item(109, 8).Value = 0
$ws.Cells.Item(109, 10).Value = 0
$ws.Cells.Item(109, 12).Value = $null
$ws.Cells.Item(109, 14).Value = 0

$ws = $wb.Worksheets.Item("CUL")  # row 38
$ws.Cells.Item(38, 8).Value = 86.333336
$ws.Cells.Item(38, 9).Value = 0
$ws.Cells.Item(38, 10).Value = 86.333336
$ws.Cells.Item(38, 11).Value = 0
$ws.Cells.Item(38, 12).Value = $null
$ws.Cells.Item(38, 13).Value = 259.000008
$ws.Cells.Item(38, 14).Value = -953.000008

$ws = $wb.Worksheets.Item("GSM")  # row 7
$ws.Cells.Item(7, 8).Value = 7000.25
$ws.Cells.Item(7, 9).Value = 5500.5
$ws.Cells.Item(7, 10).Value = 8500
$ws.Cells.Item(7, 11).Value = 5500.5
$ws.Cells.Item(7, 12).Value = 8500
$ws.Cells.Item(7, 13).Value = -5388.5
$ws.Cells.Item(7, 14).Value = -8724

$ws = $wb.Worksheets.Item("GSM")  # row 8
$ws.Cells.Item(8, 8).Value = 7000.25
$ws.Cells.Item(8, 9).Value = 5500.5
$ws.Cells.Item(8, 10).Value = 8500
$ws.Cells.Item(8, 11).Value = 5500.5
$ws.Cells.Item(8, 12).Value = 8500
$ws.Cells.Item(8, 13).Value = -5361.5
$ws.Cells.Item(8, 14).Value = -8778

$ws = $wb.Worksheets.Item("GSM")  # row 11
$ws.Cells.Item(11, 8).Value = 3680000
$ws.Cells.Item(11, 9).Value = 2977777.8
$ws.Cells.Item(11, 10).Value = 10000000
$ws.Cells.Item(11, 11).Value = 2977777.8
$ws.Cells.Item(11, 12).Value = 10000000
$ws.Cells.Item(11, 13).Value = -2977638.8
$ws.Cells.Item(11, 14).Value = -10000278

$ws = $wb.Worksheets.Item("GSM")  # row 20
$ws.Cells.Item(20, 8).Value = 0
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 0
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = $null
$ws.Cells.Item(20, 13).Value = $null
$ws.Cells.Item(20, 14).Value = 0

$ws = $wb.Worksheets.Item("GSM")  # row 92
$ws.Cells.Item(92, 8).Value = 7914
$ws.Cells.Item(92, 10).Value = 7914
$ws.Cells.Item(92, 12).Value = 7914
$ws.Cells.Item(92, 14).Value = -11658

$ws = $wb.Worksheets.Item("LTW")  # row 5
$ws.Cells.Item(5, 8).Value = 22333.334
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 10).Value = 22333.334
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 12).Value = $null
$ws.Cells.Item(5, 13).Value = 22333.334
$ws.Cells.Item(5, 14).Value = -22559.334

$ws = $wb.Worksheets.Item("LTW")  # row 20
$ws.Cells.Item(20, 8).Value = 14902.5
$ws.Cells.Item(20, 9).Value = 11005
$ws.Cells.Item(20, 11).Value = 11005
$ws.Cells.Item(20, 13).Value = -10779

$ws = $wb.Worksheets.Item("LTW")  # row 21
$ws.Cells.Item(21, 8).Value = 16971.2
$ws.Cells.Item(21, 9).Value = 11618.667
$ws.Cells.Item(21, 10).Value = 25000
$ws.Cells.Item(21, 11).Value = 11618.667
$ws.Cells.Item(21, 12).Value = 25000
$ws.Cells.Item(21, 13).Value = -11444.667
$ws.Cells.Item(21, 14).Value = -25348

$ws = $wb.Worksheets.Item("LTW")  # row 24
$ws.Cells.Item(24, 8).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 10).Value = 0
$ws.Cells.Item(24, 11).Value = 0
$ws.Cells.Item(24, 12).Value = $null
$ws.Cells.Item(24, 13).Value = $null
$ws.Cells.Item(24, 14).Value = 0

$ws = $wb.Worksheets.Item("LTW")  # row 42
$ws.Cells.Item(42, 8).Value = 100000
$ws.Cells.Item(42, 9).Value = 100000
$ws.Cells.Item(42, 11).Value = 100000
$ws.Cells.Item(42, 13).Value = -99437

$ws = $wb.Worksheets.Item("LTW")  # row 43
$ws.Cells.Item(43, 8).Value = 26933.334
$ws.Cells.Item(43, 9).Value = 22000
$ws.Cells.Item(43, 10).Value = 29400
$ws.Cells.Item(43, 11).Value = 22000
$ws.Cells.Item(43, 12).Value = 29400
$ws.Cells.Item(43, 13).Value = -21807
$ws.Cells.Item(43, 14).Value = -29786

$ws = $wb.Worksheets.Item("LTW")  # row 49
$ws.Cells.Item(49, 8).Value = 100000
$ws.Cells.Item(49, 9).Value = 100000
$ws.Cells.Item(49, 11).Value = 100000
$ws.Cells.Item(49, 13).Value = -99853

$ws = $wb.Worksheets.Item("LTW")  # row 56
$ws.Cells.Item(56, 8).Value = 14887.25
$ws.Cells.Item(56, 9).Value = 14683
$ws.Cells.Item(56, 10).Value = 15500
$ws.Cells.Item(56, 11).Value = 14683
$ws.Cells.Item(56, 12).Value = 15500
$ws.Cells.Item(56, 13).Value = -13992
$ws.Cells.Item(56, 14).Value = -16882

$ws = $wb.Worksheets.Item("WVR")  # row 3
$ws.Cells.Item(3, 8).Value = 720828
$ws.Cells.Item(3, 9).Value = 2500500
$ws.Cells.Item(3, 10).Value = 8959.200000000001
$ws.Cells.Item(3, 11).Value = 2500500
$ws.Cells.Item(3, 12).Value = 8959.200000000001
$ws.Cells.Item(3, 13).Value = -2500386
$ws.Cells.Item(3, 14).Value = -9187.200000000001

$ws = $wb.Worksheets.Item("WVR")  # row 11
$ws.Cells.Item(11, 8).Value = 5433.1665
$ws.Cells.Item(11, 9).Value = 4000
$ws.Cells.Item(11, 10).Value = 5719.8
$ws.Cells.Item(11, 11).Value = 4000
$ws.Cells.Item(11, 12).Value = 5719.8
$ws.Cells.Item(11, 13).Value = -3858
$ws.Cells.Item(11, 14).Value = -6003.8

$ws = $wb.Worksheets.Item("WVR")  # row 15
$ws.Cells.Item(15, 8).Value = 19900
$ws.Cells.Item(15, 10).Value = 19900
$ws.Cells.Item(15, 12).Value = 19900
$ws.Cells.Item(15, 14).Value = -20476

$ws = $wb.Worksheets.Item("WVR")  # row 20
$ws.Cells.Item(20, 8).Value = 22511
$ws.Cells.Item(20, 9).Value = 0
$ws.Cells.Item(20, 10).Value = 22511
$ws.Cells.Item(20, 11).Value = 0
$ws.Cells.Item(20, 12).Value = $null
$ws.Cells.Item(20, 13).Value = 22511
$ws.Cells.Item(20, 14).Value = -22991

$ws = $wb.Worksheets.Item("WVR")  # row 21
$ws.Cells.Item(21, 8).Value = 1894522
$ws.Cells.Item(21, 9).Value = 7500000
$ws.Cells.Item(21, 10).Value = 26029.334
$ws.Cells.Item(21, 11).Value = 7500000
$ws.Cells.Item(21, 12).Value = 26029.334
$ws.Cells.Item(21, 13).Value = -7499765
$ws.Cells.Item(21, 14).Value = -26499.334

$ws = $wb.Worksheets.Item("WVR")  # row 22
$ws.Cells.Item(22, 8).Value = 22222
$ws.Cells.Item(22, 9).Value = 0
$ws.Cells.Item(22, 10).Value = 22222
$ws.Cells.Item(22, 11).Value = 0
$ws.Cells.Item(22, 12).Value = $null
$ws.Cells.Item(22, 13).Value = 22222
$ws.Cells.Item(22, 14).Value = -22808

$ws = $wb.Worksheets.Item("WVR")  # row 24
$ws.Cells.Item(24, 8).Value = 1261800
$ws.Cells.Item(24, 9).Value = 2501000
$ws.Cells.Item(24, 11).Value = 2501000
$ws.Cells.Item(24, 13).Value = -2500770

$ws = $wb.Worksheets.Item("WVR")  # row 28
$ws.Cells.Item(28, 8).Value = 6999.4
$ws.Cells.Item(28, 10).Value = 4998.5
$ws.Cells.Item(28, 12).Value = 4998.5
$ws.Cells.Item(28, 14).Value = -5694.5

$ws = $wb.Worksheets.Item("WVR")  # row 30
$ws.Cells.Item(30, 8).Value = 10000
$ws.Cells.Item(30, 9).Value = 10000
$ws.Cells.Item(30, 11).Value = 10000
$ws.Cells.Item(30, 13).Value = -9893

$ws = $wb.Worksheets.Item("WVR")  # row 31
$ws.Cells.Item(31, 8).Value = 13499.5
$ws.Cells.Item(31, 9).Value = 7000
$ws.Cells.Item(31, 10).Value = 19999
$ws.Cells.Item(31, 11).Value = 7000
$ws.Cells.Item(31, 12).Value = 19999
$ws.Cells.Item(31, 13).Value = -6652
$ws.Cells.Item(31, 14).Value = -20695

$ws = $wb.Worksheets.Item("WVR")  # row 35
$ws.Cells.Item(35, 8).Value = 1894522
$ws.Cells.Item(35, 9).Value = 7500000
$ws.Cells.Item(35, 10).Value = 26029.334
$ws.Cells.Item(35, 11).Value = 7500000
$ws.Cells.Item(35, 12).Value = 26029.334
$ws.Cells.Item(35, 13).Value = -7499710
$ws.Cells.Item(35, 14).Value = -26609.334

$ws = $wb.Worksheets.Item("WVR")  # row 51
$ws.Cells.Item(51, 8).Value = 24500
$ws.Cells.Item(51, 9).Value = 24500
$ws.Cells.Item(51, 11).Value = 24500
$ws.Cells.Item(51, 13).Value = -23990

$ws = $wb.Worksheets.Item("WVR")  # row 52
$ws.Cells.Item(52, 8).Value = 10012500
$ws.Cells.Item(52, 9).Value = 20000000
$ws.Cells.Item(52, 10).Value = 25000
$ws.Cells.Item(52, 11).Value = 20000000
$ws.Cells.Item(52, 12).Value = 25000
$ws.Cells.Item(52, 13).Value = -19999774
$ws.Cells.Item(52, 14).Value = -25452

$ws = $wb.Worksheets.Item("WVR")  # row 58
$ws.Cells.Item(58, 8).Value = 0
$ws.Cells.Item(58, 9).Value = 0
$ws.Cells.Item(58, 11).Value = 0
$ws.Cells.Item(58, 13).Value = $null

$ws = $wb.Worksheets.Item("WVR")  # row 59
$ws.Cells.Item(59, 8).Value = 19000
$ws.Cells.Item(59, 10).Value = 19000
$ws.Cells.Item(59, 12).Value = 19000
$ws.Cells.Item(59, 14).Value = -20476

$ws = $wb.Worksheets.Item("WVR")  # row 122
$ws.Cells.Item(122, 8).Value = 2779.5454
$ws.Cells.Item(122, 9).Value = 2930
$ws.Cells.Item(122, 10).Value = 2102.5
$ws.Cells.Item(122, 11).Value = 8790
$ws.Cells.Item(122, 12).Value = 6307.5
$ws.Cells.Item(122, 13).Value = -6340
$ws.Cells.Item(122, 14).Value = -11207.5

$ws = $wb.Worksheets.Item("WVR")  # row 123
$ws.Cells.Item(123, 8).Value = 47499.5
$ws.Cells.Item(123, 10).Value = 47499.5
$ws.Cells.Item(123, 12).Value = 47499.5
